# Update the instructor director's name/title in the certificate template.
# "D. José Manuel Pastor Lillo, en calidad de Director" (male director)
# becomes "Dña. María Teresa Juan Díaz, en calidad de Directora" (female director).

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "D. José Manuel Pastor Lillo, en calidad de Director del I.E.S. Mare Nostrum de Alicante, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Dña. María Teresa Juan Díaz, en calidad de Directora del I.E.S. Mare Nostrum de Alicante, ",
    2
)
